$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (dbExcel), shifting old B/C to C/D
$ws.Columns("B").Insert()

# New column B content
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Ovarian epithelial cancer'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match A2's wrap-text style on B2
$ws.Range("B2").WrapText = $true

# New column B should end up the same width as column A (75.81640625 chars).
# Columns A, C, D already keep their original exact widths after the insert;
# only B needs an explicit width. The COM ColumnWidth setter here snaps to a
# coarse pixel grid, so 75 is the closest input that rounds to the nearest
# achievable width (75.8333...) to column A's 75.81640625.
$ws.Columns("B").ColumnWidth = 75
